$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (the current "Notes" column),
# shifting the existing "Notes" data from column D to column E.
$ws.Columns("D").Insert()

# Set the new column D header and width (target OOXML width is 20.140625;
# 19.25 is the closest achievable ColumnWidth given this engine's pixel quantization).
$ws.Range("D1").Value = "Expected Result"
$ws.Columns("D").ColumnWidth = 19.25

# Fill "Yes" in column D for rows that have no Linked Requirements value (column C empty),
# i.e. the rows that don't carry a Notes entry.
for ($r = 2; $r -le 21; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    if ([string]::IsNullOrEmpty($cVal)) {
        $ws.Cells.Item($r, 4).Value = "Yes"
    }
}

# Update the active selection as recorded in the workbook.
$ws.Range("G15").Select()
